$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 1375.25
$ws.Range("I2").Value = 1333.6666
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 1333.6666
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -1220.6666
$ws.Range("N2").Value = -1726

# Row 19 (Leve Item ID 7015)
$ws.Range("H19").Value = 1206.8125
$ws.Range("I19").Value = 1098.3334
$ws.Range("J19").Value = 1532.25
$ws.Range("K19").Value = 1098.3334
$ws.Range("L19").Value = 1532.25
$ws.Range("M19").Value = -923.3334
$ws.Range("N19").Value = -1882.25

# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 4299.2
$ws.Range("J43").Value = 4832.6665
$ws.Range("L43").Value = 4832.6665
$ws.Range("N43").Value = -4970.6665

# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()

# Row 111 (Leve Item ID 27768)
$ws.Range("H111").Value = 2327.1538
$ws.Range("I111").Value = 2386.7273
$ws.Range("J111").Value = 1999.5
$ws.Range("K111").Value = 7160.1819
$ws.Range("L111").Value = 5998.5
$ws.Range("M111").Value = -4093.1819
$ws.Range("N111").Value = -12132.5

# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 3033.0667
$ws.Range("I116").Value = 2177
$ws.Range("K116").Value = 2177
$ws.Range("M116").Value = 1265

# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# Row 131 (Leve Item ID 36108)
$ws.Range("H131").Value = 9779.799999999999
$ws.Range("I131").Value = 9779.799999999999
$ws.Range("K131").Value = 29339.4
$ws.Range("M131").Value = -24299.4

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 1467.2667
$ws.Range("I132").Value = 1263.4
$ws.Range("K132").Value = 3790.2
$ws.Range("M132").Value = -1260.2

# Row 141 (Leve Item ID 44161)
$ws.Range("H141").Value = 5000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 15000
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 710.4706
$ws.Range("I2").Value = 763.4167
$ws.Range("K2").Value = 763.4167
$ws.Range("M2").Value = -650.4167

# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 3092.0688
$ws.Range("I32").Value = 2617.4443
$ws.Range("K32").Value = 2617.4443
$ws.Range("M32").Value = -2330.4443

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 1920.4
$ws.Range("I61").Value = 1599.4166
$ws.Range("K61").Value = 1599.4166
$ws.Range("M61").Value = -1387.4166

# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("N110").ClearContents()

# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 710.4706
$ws.Range("I116").Value = 763.4167
$ws.Range("K116").Value = 763.4167
$ws.Range("M116").Value = 1530.5833

# Row 119 (Leve Item ID 26287)
$ws.Range("H119").Value = 120000
$ws.Range("J119").Value = 120000
$ws.Range("L119").Value = 120000
$ws.Range("N119").Value = -129676

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 1999.561
$ws.Range("I132").Value = 922.7692
$ws.Range("K132").Value = 2768.3076
$ws.Range("M132").Value = -238.3076000000001

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 1920.4
$ws.Range("I136").Value = 1599.4166
$ws.Range("K136").Value = 4798.2498
$ws.Range("M136").Value = -2248.2498

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 710.4706
$ws.Range("I3").Value = 763.4167
$ws.Range("K3").Value = 763.4167
$ws.Range("M3").Value = -649.4167

# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 3510.25
$ws.Range("I86").Value = 3315.5
$ws.Range("K86").Value = 3315.5
$ws.Range("M86").Value = -2192.5

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 3510.25
$ws.Range("I89").Value = 3315.5
$ws.Range("K89").Value = 16577.5
$ws.Range("M89").Value = -10961.5

# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 2666.5
$ws.Range("I107").Value = 2000
$ws.Range("J107").Value = 3333
$ws.Range("K107").Value = 2000
$ws.Range("L107").Value = 3333
$ws.Range("M107").Value = -80
$ws.Range("N107").Value = -7173

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 4068.5
$ws.Range("I134").Value = 4068.5
$ws.Range("K134").Value = 12205.5
$ws.Range("M134").Value = -9670.5

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 999
$ws.Range("I16").Value = 999
$ws.Range("K16").Value = 999
$ws.Range("M16").Value = -712

# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 999
$ws.Range("I113").Value = 999
$ws.Range("K113").Value = 999
$ws.Range("M113").Value = 1171

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2488.5
$ws.Range("I132").Value = 1481.5
$ws.Range("K132").Value = 4444.5
$ws.Range("M132").Value = -1914.5

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 2951
$ws.Range("I134").Value = 2951
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 8853
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -6318
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 34 (Leve Item ID 4749)
$ws.Range("H34").Value = 1975.5
$ws.Range("I34").Value = 899
$ws.Range("J34").Value = 2334.3333
$ws.Range("K34").Value = 2697
$ws.Range("L34").Value = 7002.999899999999
$ws.Range("N34").Value = -7170.999899999999
$ws.Range("M34").Value = -2613

# Row 36 (Leve Item ID 4732)
$ws.Range("H36").Value = 500
$ws.Range("J36").Value = 500
$ws.Range("L36").Value = 1500
$ws.Range("N36").Value = -1838

# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 1416.6666
$ws.Range("J131").Value = 1626.3
$ws.Range("L131").Value = 4878.9
$ws.Range("N131").Value = -14958.9

# Row 133 (Leve Item ID 44073)
$ws.Range("H133").Value = 9265
$ws.Range("I133").Value = 9265
$ws.Range("K133").Value = 27795
$ws.Range("M133").Value = -22735

$ws = $wb.Worksheets.Item("GSM")
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 3318.6667
$ws.Range("J132").Value = 3686.7144
$ws.Range("L132").Value = 11060.1432
$ws.Range("N132").Value = -16120.1432

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 2377.4546
$ws.Range("I40").Value = 2377.4546
$ws.Range("K40").Value = 2377.4546
$ws.Range("M40").Value = -2241.4546

# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 4497.375
$ws.Range("I61").Value = 4426.2856
$ws.Range("K61").Value = 4426.2856
$ws.Range("M61").Value = -4224.2856

# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 4497.375
$ws.Range("I113").Value = 4426.2856
$ws.Range("K113").Value = 4426.2856
$ws.Range("M113").Value = -2256.2856

# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 2250
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -12400

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 4384.385
$ws.Range("I132").Value = 4000.8
$ws.Range("J132").Value = 4624.125
$ws.Range("K132").Value = 12002.4
$ws.Range("L132").Value = 13872.375
$ws.Range("M132").Value = -9472.400000000001
$ws.Range("N132").Value = -18932.375

$ws = $wb.Worksheets.Item("WVR")
# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 1032.1
$ws.Range("I113").Value = 1038.6666
$ws.Range("J113").Value = 1022.25
$ws.Range("K113").Value = 3115.9998
$ws.Range("L113").Value = 3066.75
$ws.Range("M113").Value = -945.9998000000001
$ws.Range("N113").Value = -7406.75

# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 1899
$ws.Range("I122").Value = 1899
$ws.Range("K122").Value = 5697
$ws.Range("M122").Value = -3247

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 3296.5833
$ws.Range("I132").Value = 2986.75
$ws.Range("J132").Value = 3606.4167
$ws.Range("K132").Value = 8960.25
$ws.Range("L132").Value = 10819.2501
$ws.Range("M132").Value = -6430.25
$ws.Range("N132").Value = -15879.2501

